# Update of league bases, 19-04-2024 00:38
# Australia ALeague.xlsx

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Rows 104 / 105: the two match records were swapped (row 104 now holds
# what used to be row 105's data, and vice versa). Column A (row id) and
# columns C/D/E (Div / Div Original Name / Date) are identical between
# the two rows, so only B and F..AC need to move.
# ---------------------------------------------------------------------
$ws.Range("B104").Value = 7127370
$ws.Range("F104").Value = "Macarthur FC"
$ws.Range("G104").Value = "Wellington Phoenix"
$ws.Range("H104").Value = 1
$ws.Range("I104").Value = 2
$ws.Range("J104").Value = "A"
$ws.Range("K104").Value = 2.4
$ws.Range("L104").Value = 3.75
$ws.Range("M104").Value = 2.625
$ws.Range("N104").Value = 2.375
$ws.Range("O104").Value = 3.8
$ws.Range("P104").Value = 2.75
$ws.Range("Q104").Value = 0
$ws.Range("R104").Value = 1.8
$ws.Range("S104").Value = 2.05
$ws.Range("T104").Value = 3
$ws.Range("U104").Value = 1.9
$ws.Range("V104").Value = 1.95
$ws.Range("W104").Value = -1
$ws.Range("X104").Value = -1
$ws.Range("Y104").Value = 1.75
$ws.Range("Z104").Value = -1
$ws.Range("AA104").Value = 1.05
$ws.Range("AB104").Value = 0
$ws.Range("AC104").Value = 0

$ws.Range("B105").Value = 7127374
$ws.Range("F105").Value = "Central Coast Mariners"
$ws.Range("G105").Value = "Western Sydney Wanderers"
$ws.Range("H105").Value = 1
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = "H"
$ws.Range("K105").Value = 1.909
$ws.Range("L105").Value = 3.75
$ws.Range("M105").Value = 3.6
$ws.Range("N105").Value = 2.15
$ws.Range("O105").Value = 3.6
$ws.Range("P105").Value = 3.25
$ws.Range("Q105").Value = -0.25
$ws.Range("R105").Value = 1.86
$ws.Range("S105").Value = 2.04
$ws.Range("T105").Value = 2.75
$ws.Range("U105").Value = 1.975
$ws.Range("V105").Value = 1.875
$ws.Range("W105").Value = 1.15
$ws.Range("X105").Value = -1
$ws.Range("Y105").Value = -1
$ws.Range("Z105").Value = 0.8600000000000001
$ws.Range("AA105").Value = -1
$ws.Range("AB105").Value = -1
$ws.Range("AC105").Value = 0.875

# ---------------------------------------------------------------------
# Rows 124 / 125: same kind of swap as above.
# ---------------------------------------------------------------------
$ws.Range("B124").Value = 7127388
$ws.Range("F124").Value = "Sydney FC"
$ws.Range("G124").Value = "Brisbane Roar"
$ws.Range("H124").Value = 1
$ws.Range("I124").Value = 1
$ws.Range("J124").Value = "D"
$ws.Range("K124").Value = 1.5
$ws.Range("L124").Value = 5
$ws.Range("M124").Value = 5
$ws.Range("N124").Value = 1.533
$ws.Range("O124").Value = 5.25
$ws.Range("P124").Value = 5
$ws.Range("Q124").Value = -1
$ws.Range("R124").Value = 1.8
$ws.Range("S124").Value = 2.05
$ws.Range("T124").Value = 3.5
$ws.Range("U124").Value = 1.925
$ws.Range("V124").Value = 1.925
$ws.Range("W124").Value = -1
$ws.Range("X124").Value = 4.25
$ws.Range("Y124").Value = -1
$ws.Range("Z124").Value = -1
$ws.Range("AA124").Value = 1.05
$ws.Range("AB124").Value = -1
$ws.Range("AC124").Value = 0.925

$ws.Range("B125").Value = 7128012
$ws.Range("F125").Value = "Macarthur FC"
$ws.Range("G125").Value = "Central Coast Mariners"
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 3
$ws.Range("J125").Value = "A"
$ws.Range("K125").Value = 2.4
$ws.Range("L125").Value = 3.5
$ws.Range("M125").Value = 2.75
$ws.Range("N125").Value = 3.4
$ws.Range("O125").Value = 3.75
$ws.Range("P125").Value = 2.05
$ws.Range("Q125").Value = 0.25
$ws.Range("R125").Value = 2.025
$ws.Range("S125").Value = 1.825
$ws.Range("T125").Value = 3
$ws.Range("U125").Value = 2.05
$ws.Range("V125").Value = 1.8
$ws.Range("W125").Value = -1
$ws.Range("X125").Value = -1
$ws.Range("Y125").Value = 1.05
$ws.Range("Z125").Value = -1
$ws.Range("AA125").Value = 0.825
$ws.Range("AB125").Value = 0
$ws.Range("AC125").Value = 0

# ---------------------------------------------------------------------
# Row 153: closing-odds refresh.
# ---------------------------------------------------------------------
$ws.Range("N153").Value = 3
$ws.Range("O153").Value = 4
$ws.Range("P153").Value = 2.15
$ws.Range("R153").Value = 1.9
$ws.Range("S153").Value = 2
$ws.Range("U153").Value = 1.9
$ws.Range("V153").Value = 1.95

# ---------------------------------------------------------------------
# Row 154: closing-odds refresh.
# ---------------------------------------------------------------------
$ws.Range("R154").Value = 2.04
$ws.Range("S154").Value = 1.86
$ws.Range("U154").Value = 1.825
$ws.Range("V154").Value = 2.025

# ---------------------------------------------------------------------
# Row 155: closing-odds refresh.
# ---------------------------------------------------------------------
$ws.Range("U155").Value = 1.825
$ws.Range("V155").Value = 2.025

# ---------------------------------------------------------------------
# Row 156: closing-odds refresh.
# ---------------------------------------------------------------------
$ws.Range("R156").Value = 2.02
$ws.Range("S156").Value = 1.88
$ws.Range("U156").Value = 1.875
$ws.Range("V156").Value = 1.975

# ---------------------------------------------------------------------
# Row 157 takes on what used to be row 158's match (id/date/teams), with
# refreshed odds, and the old row 158 is deleted outright (the sheet
# shrinks from 158 to 157 rows).
# ---------------------------------------------------------------------
$ws.Range("B157").Value = 7127414
$ws.Range("E157").Value = 45403.16666666666
$ws.Range("F157").Value = "Perth Glory"
$ws.Range("G157").Value = "Western United FC"
$ws.Range("K157").Value = 2.4
$ws.Range("L157").Value = 3.6
$ws.Range("M157").Value = 2.625
$ws.Range("N157").Value = 2.3
$ws.Range("O157").Value = 3.75
$ws.Range("P157").Value = 2.6
$ws.Range("Q157").Value = 0
$ws.Range("R157").Value = 1.83
$ws.Range("S157").Value = 2.07
$ws.Range("T157").Value = 3.5
$ws.Range("U157").Value = 1.95
$ws.Range("V157").Value = 1.9

$ws.Rows.Item(158).Delete()
